$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "https://www.rbi.org.in"
$ws.Range("A3").Value = "https://www.owasp.org"
$ws.Range("A4").Value = "https://www.sebi.gov.in"
$ws.Range("A6").Value = "https://www.surveymonkey.com"
$ws.Range("A7").Value = "https://www.google.com"
